$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 242821.38
$ws.Range("J63").Value = 242821.38
$ws.Range("L63").Value = 242821.38
$ws.Range("N63").Value = -244069.38

$ws.Range("H66").Value = 242821.38
$ws.Range("J66").Value = 242821.38
$ws.Range("L66").Value = 728464.14
$ws.Range("N66").Value = -734704.14

$ws.Range("H75").Value = 36750
$ws.Range("J75").Value = 36750
$ws.Range("L75").Value = 36750
$ws.Range("N75").Value = -38622

$ws.Range("H78").Value = 36750
$ws.Range("J78").Value = 36750
$ws.Range("L78").Value = 110250
$ws.Range("N78").Value = -119610

$ws.Range("H123").Value = 47229.617
$ws.Range("J123").Value = 47229.617
$ws.Range("L123").Value = 47229.617
$ws.Range("N123").Value = -57029.617

$ws.Range("H124").Value = 41309
$ws.Range("J124").Value = 41309
$ws.Range("L124").Value = 41309
$ws.Range("N124").Value = -51129

$ws.Range("H126").Value = 45210
$ws.Range("J126").Value = 45210
$ws.Range("L126").Value = 45210
$ws.Range("N126").Value = -55090

$ws.Range("H134").Value = 43167.8
$ws.Range("J134").Value = 43167.8
$ws.Range("L134").Value = 43167.8
$ws.Range("N134").Value = -53307.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17407.525
$ws.Range("I32").Value = 4211.07
$ws.Range("K32").Value = 4211.07
$ws.Range("M32").Value = -3924.07

$ws.Range("H122").Value = 2104.5715
$ws.Range("I122").Value = 2052.5386
$ws.Range("K122").Value = 6157.6158
$ws.Range("M122").Value = -3707.6158

$ws.Range("H123").Value = 48500
$ws.Range("J123").Value = 48500
$ws.Range("L123").Value = 48500
$ws.Range("N123").Value = -58300

$ws.Range("H125").Value = 450027500
$ws.Range("J125").Value = 450027500
$ws.Range("L125").Value = 450027500
$ws.Range("N125").Value = -450037340

$ws.Range("H127").Value = 34650
$ws.Range("J127").Value = 34650
$ws.Range("L127").Value = 34650
$ws.Range("N127").Value = -44570

$ws.Range("H128").Value = 53444
$ws.Range("J128").Value = 53444
$ws.Range("L128").Value = 53444
$ws.Range("N128").Value = -63404

$ws.Range("H131").Value = 45509.25
$ws.Range("J131").Value = 45509.25
$ws.Range("L131").Value = 45509.25
$ws.Range("N131").Value = -55589.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 34971
$ws.Range("J126").Value = 34971
$ws.Range("L126").Value = 34971
$ws.Range("N126").Value = -44851

$ws.Range("H130").Value = 48480
$ws.Range("J130").Value = 48480
$ws.Range("L130").Value = 48480
$ws.Range("N130").Value = -58520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 44498.168
$ws.Range("J20").Value = 44498.168
$ws.Range("L20").Value = 44498.168
$ws.Range("N20").Value = -44970.168

$ws.Range("H30").Value = 44498.168
$ws.Range("J30").Value = 44498.168
$ws.Range("L30").Value = 44498.168
$ws.Range("N30").Value = -44680.168

$ws.Range("H31").Value = 1801.0952
$ws.Range("I31").Value = 1440.8055
$ws.Range("K31").Value = 1440.8055
$ws.Range("M31").Value = -1145.8055

$ws.Range("H34").Value = 1801.0952
$ws.Range("I34").Value = 1440.8055
$ws.Range("K34").Value = 1440.8055
$ws.Range("M34").Value = -1238.8055

$ws.Range("H124").Value = 14315.546
$ws.Range("J124").Value = 14315.546
$ws.Range("L124").Value = 14315.546
$ws.Range("N124").Value = -19225.546

$ws.Range("H127").Value = 48805.715
$ws.Range("J127").Value = 48805.715
$ws.Range("L127").Value = 48805.715
$ws.Range("N127").Value = -58725.715

$ws.Range("H128").Value = 44498.168
$ws.Range("J128").Value = 44498.168
$ws.Range("L128").Value = 44498.168
$ws.Range("N128").Value = -54458.168

$ws.Range("H130").Value = 37397.145
$ws.Range("J130").Value = 37397.145
$ws.Range("L130").Value = 37397.145
$ws.Range("N130").Value = -47437.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 215350.22
$ws.Range("I5").Value = 247.78572
$ws.Range("J5").Value = 2222973
$ws.Range("K5").Value = 743.35716
$ws.Range("L5").Value = 6668919
$ws.Range("M5").Value = -631.35716
$ws.Range("N5").Value = -6669143

$ws.Range("H122").Value = 50320
$ws.Range("I122").Value = 343.14285
$ws.Range("J122").Value = 60174.59
$ws.Range("K122").Value = 3088.28565
$ws.Range("L122").Value = 541571.3099999999
$ws.Range("M122").Value = -638.2856500000003
$ws.Range("N122").Value = -546471.3099999999

$ws.Range("H135").Value = 215350.22
$ws.Range("I135").Value = 247.78572
$ws.Range("J135").Value = 2222973
$ws.Range("K135").Value = 2230.07148
$ws.Range("L135").Value = 20006757
$ws.Range("M135").Value = 304.9285199999999
$ws.Range("N135").Value = -20011827

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2458.3333
$ws.Range("I122").Value = 2062.5
$ws.Range("K122").Value = 6187.5
$ws.Range("M122").Value = -3737.5

$ws.Range("H128").Value = 45411.11
$ws.Range("J128").Value = 45411.11
$ws.Range("L128").Value = 45411.11
$ws.Range("N128").Value = -55371.11

$ws.Range("H132").Value = 1721.0172
$ws.Range("I132").Value = 1425.85
$ws.Range("K132").Value = 4277.549999999999
$ws.Range("M132").Value = -1747.549999999999

$ws.Range("H135").Value = 60930.77
$ws.Range("J135").Value = 60930.77
$ws.Range("L135").Value = 60930.77
$ws.Range("N135").Value = -71070.76999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3720.4348
$ws.Range("I122").Value = 2386.25
$ws.Range("J122").Value = 4432
$ws.Range("K122").Value = 7158.75
$ws.Range("L122").Value = 13296
$ws.Range("M122").Value = -4708.75
$ws.Range("N122").Value = -18196

$ws.Range("H128").Value = 28569.857
$ws.Range("J128").Value = 28569.857
$ws.Range("L128").Value = 28569.857
$ws.Range("N128").Value = -38529.857

$ws.Range("H130").Value = 48457.5
$ws.Range("J130").Value = 48457.5
$ws.Range("L130").Value = 48457.5
$ws.Range("N130").Value = -58497.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

$ws.Range("H127").Value = 48768.125
$ws.Range("J127").Value = 48768.125
$ws.Range("L127").Value = 48768.125
$ws.Range("N127").Value = -58688.125

$ws.Range("H135").Value = 57664
$ws.Range("J135").Value = 57664
$ws.Range("L135").Value = 57664
$ws.Range("N135").Value = -67804
